# Apply the LinuxForHealth rebrand + version/date bump to the
# StructureDefinition-insured workbook.
#
# The workbook has two sheets:
#   "Metadata" - simple Property/Value table
#   "Elements" - wide FHIR element grid
#
# Only the visible cell contents below need to change; Excel will
# rebuild its own shared-string table on save.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insured"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
# Root "Extension" row no longer repeats the full constraint expression
# in the Constraint(s) column.
$elements.Range("AI2").Value = ""

# insuredRank / insuredCategory extension type URLs + the url fixed value
# all move from ibm.com to linuxforhealth.org.
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insured-rank}
"
$elements.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insured-category}
"
$elements.Range("Q7").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insured"
